# [IMP] Change column header in template
#
# Row 13 holds the report's column headers. Three header labels change:
#   D13 "Responsible" -> "Validated By"
#   H13 "Amount"       -> "PV Amount"
#   R13 "Amount"       -> "KV Amount"   (was sharing text with H13, now distinct)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D13").Value = "Validated By"
$ws.Range("H13").Value = "PV Amount"
$ws.Range("R13").Value = "KV Amount"
